# Overview_demos.xlsx - "updates figs and demos"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Highlight the "No figures yet" remarks that still need attention
#    (D20, D22, D24, D25, D26) with a themed fill (Green, Accent 6).
# ---------------------------------------------------------------------
$ws.Range("D20,D22,D24,D25,D26").Interior.ThemeColor = 10

# ---------------------------------------------------------------------
# 2. Append the newly-translated Dutch demos (rows 27-36).
#    Values are entered in the exact order the original author typed
#    them so the shared-string table grows in the same sequence.
# ---------------------------------------------------------------------
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "1 01"
$ws.Range("C27").Value = "Opwaartse en Neerwaartsekracht"

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "1 03"
$ws.Range("C28").Value = "Waardoor stijgt het water?"

$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "1 10"
$ws.Range("C29").Value = "Warm & Koud, Zoet & Zout"

$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "1 17"
$ws.Range("C30").Value = "Blussen zonder water"

$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "1 21"
$ws.Range("C31").Value = "Lucht is niet niks"

$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "1 29"
$ws.Range("C32").Value = "Gloeilamp uitblazen"

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "1 37"
$ws.Range("C33").Value = "Twee veren, serie of parallel?"

$ws.Range("A34").Value = 33
$ws.Range("C34").Value = "Schaduw van een vlam"
$ws.Range("B34").Value = "1 52"

$ws.Range("A35").Value = 34
$ws.Range("B35").Value = "1 57"
$ws.Range("C35").Value = "Spelen met dichtheid"

$ws.Range("A36").Value = 35
$ws.Range("C36").Value = "Vallend kaars"
$ws.Range("B36").Value = "1 59"

$ws.Range("D27").Value = "NL"
$ws.Range("D28").Value = "NL"
$ws.Range("D29").Value = "NL"
$ws.Range("D30").Value = "NL"
$ws.Range("D31").Value = "NL"
$ws.Range("D32").Value = "NL"
$ws.Range("D33").Value = "NL"
$ws.Range("D34").Value = "NL"
$ws.Range("D35").Value = "NL"
$ws.Range("D36").Value = "NL"

# ---------------------------------------------------------------------
# 3. Extend the running demo-number count down to row 100 (values 36-99
#    in column A), reserving space for more demos to be catalogued.
# ---------------------------------------------------------------------
for ($r = 37; $r -le 100; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# ---------------------------------------------------------------------
# 4. View tidy-up: zoom out a bit and leave the selection on the newly
#    highlighted D20 cell; drop the old frozen "topLeftCell" scroll.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 85
$ws.Range("D20").Select()
